$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("O3").Value = 1.21
$ws.Range("S3").Value = 2.58
$ws.Range("T3").Value = 1.63
$ws.Range("U3").Value = 2.48

# Row 4
$ws.Range("X4").Value = 25

# Row 6
$ws.Range("L6").Value = 1.24
$ws.Range("U6").Value = 2.2
$ws.Range("AL6").Value = 26

# Row 7
$ws.Range("AN7").Value = 27

# Row 8
$ws.Range("H8").Value = 6
$ws.Range("P8").Value = 1.93
$ws.Range("U8").Value = 1.93
$ws.Range("X8").Value = 13
$ws.Range("AC8").Value = 8.800000000000001

# Row 12
$ws.Range("F12").Value = 1.99
$ws.Range("G12").Value = 2.12
$ws.Range("L12").Value = 1.01
$ws.Range("M12").Value = 1.01
$ws.Range("N12").Value = 1.94
$ws.Range("O12").Value = 1.29
$ws.Range("R12").Value = 1.31
$ws.Range("S12").Value = 2.86
$ws.Range("T12").Value = 1.64
$ws.Range("U12").Value = 1.92
$ws.Range("V12").Value = 1.29
$ws.Range("W12").Value = 1.9
$ws.Range("X12").Value = 20
$ws.Range("Y12").Value = 21
$ws.Range("Z12").Value = 44
$ws.Range("AA12").Value = 1000
$ws.Range("AB12").Value = 13
$ws.Range("AC12").Value = 12
$ws.Range("AD12").Value = 24
$ws.Range("AE12").Value = 75
$ws.Range("AF12").Value = 18
$ws.Range("AG12").Value = 15
$ws.Range("AH12").Value = 25
$ws.Range("AI12").Value = 80
$ws.Range("AJ12").Value = 34
$ws.Range("AK12").Value = 30
$ws.Range("AL12").Value = 50
$ws.Range("AM12").Value = 1000
$ws.Range("AN12").Value = 1000
$ws.Range("AO12").Value = 1000
